# Adds the Day-1 / kickup / Initial usages / Day-2 / Initial usage contd.. /
# Day-3 slides (slides 2-7) to the deck, per the "Prompting, Passing, Command
# line executions are added" commit.

$p = $ppt.ActivePresentation

$titleOnlyLayout    = $p.SlideMaster.CustomLayouts.Item(6)   # "Title Only"
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# --- Slide 2: "Day-1" --------------------------------------------------
$s2 = $p.Slides.AddSlide(2, $titleOnlyLayout)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Day-1"

# --- Slide 3: "kickup" --------------------------------------------------
$s3 = $p.Slides.AddSlide(3, $titleContentLayout)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "kickup"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Installation & setup`rAnaconda`rSpyder`rJupyter`rHelloworld`r"

# --- Slide 4: "Initial usages" ------------------------------------------
$s4 = $p.Slides.AddSlide(4, $titleContentLayout)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Initial usages"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "Comments and Pound character`rNumbers and Math`rVariables and Names`rVariables and Printing`r"

# --- Slide 5: "Day-2" ----------------------------------------------------
$s5 = $p.Slides.AddSlide(5, $titleOnlyLayout)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Day-2"

# --- Slide 6: "Initial usage contd.." ------------------------------------
$s6 = $p.Slides.AddSlide(6, $titleContentLayout)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Initial usage contd.."
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "strings and text`rmore printing`rprinting, printing`rprinting, printing, printing`rescape sequences`rasking questions --- bots`r"

# --- Slide 7: "Day-3" ------------------------------------------------------
$s7 = $p.Slides.AddSlide(7, $titleContentLayout)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Day-3"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = ""
